# Weekly refresh of the "Fruta, Vega Central Mapocho de Santiago - Coco" sheet.
# A new weekly observation is inserted as row 7 (pushing the existing rows
# 7-38 down to 8-39), growing the used range from A1:T38 to A1:T39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 7; existing rows 7..38 shift to 8..39,
# carrying their formatting (incl. the date number-format on column D) down
# with them automatically.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the new weekly record.
$ws.Range("A7").Value = 9
$ws.Range("B7").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C7").Value = "Metropolitana"
$ws.Range("D7").Value = 44620
$ws.Range("E7").Value = 13
$ws.Range("F7").Value = "Fruta"
$ws.Range("G7").Value = 100108
$ws.Range("H7").Value = "Tropicales y subtropicales"
$ws.Range("I7").Value = 100108007
$ws.Range("J7").Value = "Coco"
$ws.Range("K7").Value = "Sin especificar"
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 60
$ws.Range("N7").Value = 22000
$ws.Range("O7").Value = 22000
$ws.Range("P7").Value = 22000
$ws.Range("Q7").Value = "$/malla 20 unidades"
$ws.Range("R7").Value = "Perú"
$ws.Range("S7").Value = 1100
$ws.Range("T7").Value = 20
